$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phase 1: fill B, C, D, F columns for rows 7-20 in entry order
# Row 7
$ws.Range("F7").Value = "1:01pm"
# Row 8
$ws.Range("B8").Value = "Design"
$ws.Range("C8").Value = 43558
$ws.Range("C8").NumberFormat = "m/d/yy"
$ws.Range("D8").Value = "2:13pm"
$ws.Range("F8").Value = "2:21pm"
# Row 10
$ws.Range("B10").Value = "Design"
$ws.Range("C10").Value = 43558
$ws.Range("C10").NumberFormat = "m/d/yy"
$ws.Range("D10").Value = "2:37pm"
$ws.Range("F10").Value = "2:40pm"
# Row 20
$ws.Range("B20").Value = "Design"
$ws.Range("C20").Value = 43559
$ws.Range("C20").NumberFormat = "m/d/yy"
$ws.Range("D20").Value = "11:02am"
$ws.Range("F20").Value = "11:07am"
# Row 9
$ws.Range("B9").Value = "Testing"
$ws.Range("C9").Value = 43558
$ws.Range("C9").NumberFormat = "m/d/yy"
$ws.Range("D9").Value = "2:23pm"
$ws.Range("F9").Value = "2:32pm"
# Row 11
$ws.Range("B11").Value = "Coding"
$ws.Range("C11").Value = 43558
$ws.Range("C11").NumberFormat = "m/d/yy"
$ws.Range("D11").Value = "3:16pm"
$ws.Range("F11").Value = "3:21pm"
# Row 12
$ws.Range("B12").Value = "Coding"
$ws.Range("C12").Value = 43558
$ws.Range("C12").NumberFormat = "m/d/yy"
$ws.Range("D12").Value = "3:31pm"
$ws.Range("F12").Value = "3:35pm"
# Row 13
$ws.Range("B13").Value = "Coding"
$ws.Range("C13").Value = 43558
$ws.Range("C13").NumberFormat = "m/d/yy"
$ws.Range("D13").Value = "3:42pm"
$ws.Range("F13").Value = "3:48pm"
# Row 14
$ws.Range("B14").Value = "Coding"
$ws.Range("C14").Value = 43558
$ws.Range("C14").NumberFormat = "m/d/yy"
$ws.Range("D14").Value = "4:22pm"
$ws.Range("F14").Value = "4:31pm"
# Row 15
$ws.Range("B15").Value = "Coding"
$ws.Range("C15").Value = 43558
$ws.Range("C15").NumberFormat = "m/d/yy"
$ws.Range("D15").Value = "4:36pm"
$ws.Range("F15").Value = "4:53pm"
# Row 16
$ws.Range("B16").Value = "Coding"
$ws.Range("C16").Value = 43558
$ws.Range("C16").NumberFormat = "m/d/yy"
$ws.Range("D16").Value = "6:55pm"
$ws.Range("F16").Value = "7:05pm"
# Row 18
$ws.Range("B18").Value = "Coding"
$ws.Range("C18").Value = 43559
$ws.Range("C18").NumberFormat = "m/d/yy"
$ws.Range("D18").Value = "10:46am"
$ws.Range("F18").Value = "10:56am"
# Row 17
$ws.Range("B17").Value = "Testing"
$ws.Range("C17").Value = 43559
$ws.Range("C17").NumberFormat = "m/d/yy"
$ws.Range("D17").Value = "10:13am"
$ws.Range("F17").Value = "10:19am"
# Row 19
$ws.Range("B19").Value = "Testing"
$ws.Range("C19").Value = 43559
$ws.Range("C19").NumberFormat = "m/d/yy"
$ws.Range("D19").Value = "10:56am"
$ws.Range("F19").Value = "10:57am"

# Phase 2: fill H column for rows 8-20 in ascending order
$ws.Range("H8").Value = "Preliminary designs"
$ws.Range("H9").Value = "Write tests"
$ws.Range("H10").Value = "More design"
$ws.Range("H11").Value = "Code "
$ws.Range("H12").Value = "Code "
$ws.Range("H13").Value = "Code "
$ws.Range("H14").Value = "Code "
$ws.Range("H15").Value = "Code "
$ws.Range("H16").Value = "Code "
$ws.Range("H17").Value = "Test code"
$ws.Range("H18").Value = "Fix code"
$ws.Range("H19").Value = "Test again"
$ws.Range("H20").Value = "Fix designs"

$ws.Range("C20").Select() | Out-Null
